$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q1" right after "2021-Q4"
#    (i.e. right before "总计") and populate it with fund data.
# ---------------------------------------------------------------
$wsAfter = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Add($null, $wsAfter)
$ws.Name = "2022-Q1"

# Copy header formatting (bold + border + centered) from an existing
# quarter sheet so the new header row matches the workbook's style.
$srcHeader = $wb.Worksheets.Item("2021-Q4").Range("B1:H1")
$srcHeader.Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Copy the row-index cell formatting (bold + border, centered) used in
# column A of the existing quarter sheets.
$srcIndex = $wb.Worksheets.Item("2021-Q4").Range("A2")
$srcIndex.Copy()
$ws.Range("A2:A28").PasteSpecial(-4122)

# Columns B-G hold text values (fund code / name / numeric-looking
# strings) - force text storage before writing so values such as
# "000055" or "75.36" are not coerced into numbers.
$ws.Range("B2:G28").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "000055"
$ws.Range("C2").Value = "广发纳斯达克100指数(QDII) A 美元现汇"
$ws.Range("D2").Value = "75.36"
$ws.Range("E2").Value = "85.84"
$ws.Range("F2").Value = "6.92"
$ws.Range("G2").Value = "5.2149"
$ws.Range("H2").Value = 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "270042"
$ws.Range("C3").Value = "广发纳斯达克100指数QDII A"
$ws.Range("D3").Value = "75.36"
$ws.Range("E3").Value = "85.84"
$ws.Range("F3").Value = "6.92"
$ws.Range("G3").Value = "5.2149"
$ws.Range("H3").Value = 3
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "006479"
$ws.Range("C4").Value = "广发纳斯达克100指数（QDII）C人民币"
$ws.Range("D4").Value = "75.36"
$ws.Range("E4").Value = "85.84"
$ws.Range("F4").Value = "6.92"
$ws.Range("G4").Value = "5.2149"
$ws.Range("H4").Value = 3
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "006480"
$ws.Range("C5").Value = "广发纳斯达克100指数（QDII）C美元现汇"
$ws.Range("D5").Value = "75.36"
$ws.Range("E5").Value = "85.84"
$ws.Range("F5").Value = "6.92"
$ws.Range("G5").Value = "5.2149"
$ws.Range("H5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "513100"
$ws.Range("C6").Value = "国泰纳斯达克100 (QDII-ETF)"
$ws.Range("D6").Value = "41.86"
$ws.Range("E6").Value = "90.74"
$ws.Range("F6").Value = "6.59"
$ws.Range("G6").Value = "2.7586"
$ws.Range("H6").Value = 3
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "513500"
$ws.Range("C7").Value = "博时标普500ETF(QDII)"
$ws.Range("D7").Value = "70.03"
$ws.Range("E7").Value = "90.45"
$ws.Range("F7").Value = "3.37"
$ws.Range("G7").Value = "2.3600"
$ws.Range("H7").Value = 3
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "040046"
$ws.Range("C8").Value = "华安纳斯达克100指数QDII - 人民币"
$ws.Range("D8").Value = "22.85"
$ws.Range("E8").Value = "90.93"
$ws.Range("F8").Value = "6.61"
$ws.Range("G8").Value = "1.5104"
$ws.Range("H8").Value = 3
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "040047"
$ws.Range("C9").Value = "华安纳斯达克100指数QDII - 美元现钞"
$ws.Range("D9").Value = "22.85"
$ws.Range("E9").Value = "90.93"
$ws.Range("F9").Value = "6.61"
$ws.Range("G9").Value = "1.5104"
$ws.Range("H9").Value = 3
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "040048"
$ws.Range("C10").Value = "华安纳斯达克100指数QDII - 美元现汇"
$ws.Range("D10").Value = "22.85"
$ws.Range("E10").Value = "90.93"
$ws.Range("F10").Value = "6.61"
$ws.Range("G10").Value = "1.5104"
$ws.Range("H10").Value = 3
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "000988"
$ws.Range("C11").Value = "嘉实全球互联网股票 - 人民币QDII"
$ws.Range("D11").Value = "13.21"
$ws.Range("E11").Value = "85.88"
$ws.Range("F11").Value = "10.18"
$ws.Range("G11").Value = "1.3448"
$ws.Range("H11").Value = 3
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "000989"
$ws.Range("C12").Value = "嘉实全球互联网股票 - 美元现汇QDII"
$ws.Range("D12").Value = "13.21"
$ws.Range("E12").Value = "85.88"
$ws.Range("F12").Value = "10.18"
$ws.Range("G12").Value = "1.3448"
$ws.Range("H12").Value = 3
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "000990"
$ws.Range("C13").Value = "嘉实全球互联网股票 - 美元现钞QDII"
$ws.Range("D13").Value = "13.21"
$ws.Range("E13").Value = "85.88"
$ws.Range("F13").Value = "10.18"
$ws.Range("G13").Value = "1.3448"
$ws.Range("H13").Value = 3
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "160213"
$ws.Range("C14").Value = "国泰纳斯达克100指数(QDII)"
$ws.Range("D14").Value = "15.88"
$ws.Range("E14").Value = "90.49"
$ws.Range("F14").Value = "6.58"
$ws.Range("G14").Value = "1.0449"
$ws.Range("H14").Value = 3
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "000834"
$ws.Range("C15").Value = "大成纳斯达克100指数 (QDII)"
$ws.Range("D15").Value = "14.15"
$ws.Range("E15").Value = "89.27"
$ws.Range("F15").Value = "6.50"
$ws.Range("G15").Value = "0.9198"
$ws.Range("H15").Value = 3
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "270023"
$ws.Range("C16").Value = "广发全球精选股票(QDII)"
$ws.Range("D16").Value = "25.53"
$ws.Range("E16").Value = "78.43"
$ws.Range("F16").Value = "3.41"
$ws.Range("G16").Value = "0.8706"
$ws.Range("H16").Value = 10
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "000906"
$ws.Range("C17").Value = "广发全球精选股票(QDII)美元现汇"
$ws.Range("D17").Value = "25.53"
$ws.Range("E17").Value = "78.43"
$ws.Range("F17").Value = "3.41"
$ws.Range("G17").Value = "0.8706"
$ws.Range("H17").Value = 10
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "159941"
$ws.Range("C18").Value = "广发纳斯达克100ETFQDII"
$ws.Range("D18").Value = "11.87"
$ws.Range("E18").Value = "90.26"
$ws.Range("F18").Value = "6.89"
$ws.Range("G18").Value = "0.8178"
$ws.Range("H18").Value = 3
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "002423"
$ws.Range("C19").Value = "华宝兴业标普美国消费(QDII-LOF)美元"
$ws.Range("D19").Value = "3.62"
$ws.Range("E19").Value = "94.37"
$ws.Range("F19").Value = "21.87"
$ws.Range("G19").Value = "0.7917"
$ws.Range("H19").Value = 1
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "162415"
$ws.Range("C20").Value = "华宝标普美国消费(QDII-LOF)人民币A"
$ws.Range("D20").Value = "3.62"
$ws.Range("E20").Value = "94.37"
$ws.Range("F20").Value = "21.87"
$ws.Range("G20").Value = "0.7917"
$ws.Range("H20").Value = 1
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "000043"
$ws.Range("C21").Value = "嘉实美国成长股票(QDII) -人民币"
$ws.Range("D21").Value = "14.64"
$ws.Range("E21").Value = "94.24"
$ws.Range("F21").Value = "4.72"
$ws.Range("G21").Value = "0.6910"
$ws.Range("H21").Value = 4
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "000044"
$ws.Range("C22").Value = "嘉实美国成长股票(QDII) - 美元现汇"
$ws.Range("D22").Value = "14.64"
$ws.Range("E22").Value = "94.24"
$ws.Range("F22").Value = "4.72"
$ws.Range("G22").Value = "0.6910"
$ws.Range("H22").Value = 4
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "003722"
$ws.Range("C23").Value = "易方达纳斯达克100指数美元（QDII-LOF）"
$ws.Range("D23").Value = "9.07"
$ws.Range("E23").Value = "91.29"
$ws.Range("F23").Value = "6.66"
$ws.Range("G23").Value = "0.6041"
$ws.Range("H23").Value = 3
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "161130"
$ws.Range("C24").Value = "易方达纳斯达克100指数人民币（QDII-LOF）"
$ws.Range("D24").Value = "9.07"
$ws.Range("E24").Value = "91.29"
$ws.Range("F24").Value = "6.66"
$ws.Range("G24").Value = "0.6041"
$ws.Range("H24").Value = 3
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "003718"
$ws.Range("C25").Value = "易方达标普500指数(QDII-LOF) 美元"
$ws.Range("D25").Value = "5.22"
$ws.Range("E25").Value = "91.11"
$ws.Range("F25").Value = "3.37"
$ws.Range("G25").Value = "0.1759"
$ws.Range("H25").Value = 3
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "161125"
$ws.Range("C26").Value = "易方达标普500指数(QDII-LOF) 人民币"
$ws.Range("D26").Value = "5.22"
$ws.Range("E26").Value = "91.11"
$ws.Range("F26").Value = "3.37"
$ws.Range("G26").Value = "0.1759"
$ws.Range("H26").Value = 3
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "009975"
$ws.Range("C27").Value = "华宝标普美国消费(QDII-LOF)人民币C"
$ws.Range("D27").Value = "0.61"
$ws.Range("E27").Value = "94.37"
$ws.Range("F27").Value = "21.87"
$ws.Range("G27").Value = "0.1334"
$ws.Range("H27").Value = 1
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "006555"
$ws.Range("C28").Value = "浦银安盛全球智能科技股票（QDII）"
$ws.Range("D28").Value = "3.20"
$ws.Range("E28").Value = "85.41"
$ws.Range("F28").Value = "2.67"
$ws.Range("G28").Value = "0.0854"
$ws.Range("H28").Value = 7

# Reset the index column's value (PasteSpecial also copied the source
# cell's value of 0 into every row - Value is overwritten individually
# below for each row, so nothing further is required here).

# After writing all values, clear the temporary text-number-format so
# the cells keep their text type but drop the extra style id.
$ws.Range("B2:G28").Style = "Normal"

# Give the new sheet a correct dimension box.
$ws.Range("A1").Select()

# ---------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert a new first data row
#    for 2022-Q1 and bump the existing running index column by one.
# ---------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert()

$srcA = $tot.Range("A3")
$srcA.Copy()
$tot.Range("A2").PasteSpecial(-4122)
$tot.Range("A2").Value = 0

$tot.Range("B2").Style = "Normal"
$tot.Range("C2").Style = "Normal"
$tot.Range("D2").Style = "Normal"
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 27
$tot.Range("D2").Value = 43.81

$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4
$tot.Range("A7").Value = 5

Write-Host "edit complete"
